$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ridership")

# Update the Monday (11 Jul 2016) row values: Riders 183 -> 196, Average 86.2 -> 86.53
$ws.Range("C2").Value = 196
$ws.Range("D2").Value = 86.53

$wb.Save()
